$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "70.016.74"
$ws.Range("E2").Value = "  +0.72%  "
$ws.Range("D3").Value = "3.695.12"
$ws.Range("E3").Value = "  +0.12%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.13%  "
$ws.Range("E5").Value = "  -4.03%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "162.20"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.39%  "
$ws.Range("E7").Value = "  -0.04%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.147"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.15%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "7.21"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.07%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.448"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.85%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0000235"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.58%  "
$ws.Range("D13").Value = "4.312.86"
$ws.Range("E13").Value = "  +0.01%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "32.98"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.55%  "
$ws.Range("D15").Value = "3.697.64"
$ws.Range("E15").Value = "  -0.41%  "
$ws.Range("D16").Value = "69.933.46"
$ws.Range("E16").Value = "  +0.68%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "16.20"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.16%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.55"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.22%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "10.51"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +7.29%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "475.12"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.01%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.654"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.53%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "80.04"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.54%  "
$ws.Range("D24").Value = "3.841.23"
$ws.Range("E24").Value = "  +0.10%  "
$ws.Range("B25").Value = "PEPE"
$ws.Range("C25").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.0000128"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.28%  "
$ws.Range("B26").Value = "Dai"
$ws.Range("C26").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.00"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.04%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "11.06"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.86%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.25"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.46%  "
$ws.Range("E29").Value = "  -1.40%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.74"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.02%  "
$ws.Range("E31").Value = "  +0.65%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.60"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.23%  "
$ws.Range("B33").Value = "EthereumClassic"
$ws.Range("C33").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "27.01"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.06%  "
$ws.Range("B34").Value = "Binance-PegBSC-USD"
$ws.Range("C34").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.999"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.49%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.166"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +2.04%  "
$ws.Range("D36").Value = "3.691.86"
$ws.Range("E36").Value = "  +0.32%  "
$ws.Range("E37").Value = "  +0.22%  "
$ws.Range("E38").Value = "  -0.10%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.29"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.27%  "
$ws.Range("B40").Value = "Filecoin"
$ws.Range("C40").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.93"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -4.79%  "
$ws.Range("B41").Value = "Monero"
$ws.Range("C41").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "180.23"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +6.92%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.999"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.12%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0909"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.96%  "
$ws.Range("E44").Value = "  -1.10%  "
$ws.Range("B45").Value = "OKB"
$ws.Range("C45").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "47.13"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.96%  "
$ws.Range("B46").Value = "InjectiveProtocol"
$ws.Range("C46").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "29.31"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +5.69%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.81"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +2.64%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.000273"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.75%  "
$ws.Range("B49").Value = "ONDO"
$ws.Range("C49").Value = "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.26"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -2.84%  "
$ws.Range("B50").Value = "SuiNetwork"
$ws.Range("C50").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.07"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.65%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.89"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.00%  "
